# Adding limits on investment: a new "permit and construction" column
# (expectedLeadtime + expectedPermittime) on the TechnologiesEmlab sheet,
# plus restoring the sheet selections / active tab that were in effect
# when the workbook was last saved.

$wb = $excel.ActiveWorkbook

# --- TechnologiesEmlab: new column I = permit + construction time ---
$wsTech = $wb.Worksheets.Item("TechnologiesEmlab")
$wsTech.Range("I1").Value = "permit and construction"
for ($r = 2; $r -le 25; $r++) {
    $wsTech.Cells.Item($r, 9).Formula = "=C$r+B$r"
}
$wsTech.Range("I15").Select()

# --- Fuels: remembered selection moved to K18 ---
$wsFuels = $wb.Worksheets.Item("Fuels")
$wsFuels.Range("K18").Select()

# --- EnergyProducers: becomes the active sheet, selection moved to M15 ---
$wsProducers = $wb.Worksheets.Item("EnergyProducers")
$wsProducers.Activate()
$wsProducers.Range("M15").Select()
